$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Paragraph "Each Results_xxx...": merge "run [on] all 13 regions"
#    (text unchanged, just collapses proofErr-wrapped runs)
# -----------------------------------------------------------------
$d.Content.Find.Execute("That says it was run on all 13 regions.", $true, $false, $false, $false, $false, $true, 1, $false, "That says it was run on all 13 regions.", 2) | Out-Null

# -----------------------------------------------------------------
# 2. "Right-click on each .xlsm file, and set properties to allow macros."
#    -> "Right-click on each .xlsm file in the file manager, and set properties to allow macros."
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(10)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Right-click on each .xlsm file in the file manager, and set properties to allow macros."

# -----------------------------------------------------------------
# 3. "Open Compare.xlsm..." paragraph - big rewrite
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Open " + [char]8220 + "Compare.xlsm" + [char]8221 + ".  Enable macros here.  This should open all the Results" + [char]8230 + "xlsm files. If that did not happen, click on cell A1.  There will be an arrow to the right of the cell.  Click on this dropdown, and select something else.  Try again with A2, A3, the " + [char]8220 + "To Bar" + [char]8221 + " button.  Once they have all loaded, click on each Results_xxx.xlsm file, and enable macros in each one.  Go back to Compare.xlsm."

# -----------------------------------------------------------------
# 4. "On sheet Graph..." paragraph - big rewrite
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(12)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "On sheet Graph, cell A1, is the region.  If you click on it, an arrow appears to the right.  Hit this arrow to choose a region or Sum_all to look at entire US.  EIA divides the country into 13 regions, like California and Mid West.  "

# -----------------------------------------------------------------
# 5. "Cell A2 is what you are graphing..." - add " etc" before period
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Cell A2 is what you are graphing - Cost, Energy (MWh), MTons CO2 etc. "

# -----------------------------------------------------------------
# 6. "At the base of the model..." - append new run "  We call this year 0."
# -----------------------------------------------------------------
$p = $d.Paragraphs.Item(16)
$r = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$r.InsertBefore("  We call this year 0.")

# -----------------------------------------------------------------
# 7. "Initial demand is set at the sum of energies for each hour.  Again, 35,064 values."
#    -> insert " in year 0" before period
# -----------------------------------------------------------------
$d.Content.Find.Execute("for each hour.  Again, 35,064 values.", $true, $false, $false, $false, $false, $true, 1, $false, "for each hour in year 0.  Again, 35,064 values.", 2) | Out-Null

# -----------------------------------------------------------------
# 8. "lastRenderedPageBreak" moves from paragraph 25 to paragraph 24
# -----------------------------------------------------------------
$p25 = $d.Paragraphs.Item(25)
$p25.Range.Find.Execute("Running is very processor intensive", $true, $false, $false, $false, $false, $true, 1, $false, "Running is very processor intensive", 2) | Out-Null

Write-Host "Done"
